# The edit rotates the data held in rows 29-35 of the "Artfynd" sheet among
# themselves (same 7 observation records, reassigned to different row
# positions) while leaving every other row untouched.
#
# Mapping of "new row -> old row that its data comes from":
#   29 <- 31
#   30 <- 35
#   31 <- 32
#   32 <- 29
#   33 <- 30
#   34 <- 33
#   35 <- 34

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 29
$lastRow = 35
$firstCol = 1
$lastCol = 51

# Columns (1-based) whose values in this block are stored as plain text,
# including some that look like numbers/dates/times ("200", "2023-08-03",
# "00:00", ...). Forcing a text number format on these before the write
# keeps such values from being auto-converted into numeric / date serials.
$textCols = @(3,4,6,7,8,9,10,11,12,13,14,16,20,21,22,23,25,26,27,28,32,46,49,50,51)
$textColSet = @{}
foreach ($tc in $textCols) { $textColSet[$tc] = $true }

# Snapshot every cell value in the affected rows before writing anything
# back, since several rows are both a source and a destination for the
# rotation.
$snapshotValue = @{}
foreach ($r in $firstRow..$lastRow) {
    $rowValues = @{}
    foreach ($c in $firstCol..$lastCol) {
        $rowValues[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshotValue[$r] = $rowValues
}

$sourceForRow = @{
    29 = 31
    30 = 35
    31 = 32
    32 = 29
    33 = 30
    34 = 33
    35 = 34
}

foreach ($destRow in $firstRow..$lastRow) {
    $srcRow = $sourceForRow[$destRow]
    $rowValues = $snapshotValue[$srcRow]
    foreach ($c in $firstCol..$lastCol) {
        $cell = $ws.Cells.Item($destRow, $c)
        if ($textColSet.ContainsKey($c)) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $rowValues[$c]
    }
}

Write-Output "Rotated rows 29-35 data successfully"
